# Generate Report for Handoff
#
# File "b.md" has been handed off again: status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", and a new
# handoff file / handoff datetime is recorded for both the zh-cn and
# de-de locales.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusReady
$overview.Range("C3").Value = $statusReady
$overview.Range("D3").Value = "2016-47-09 09:47:22"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-09 09:47:11"

foreach ($hl in $zhcn.Hyperlinks) {
    $r = $hl.Range
    if ($r.Row -eq 3 -and $r.Column -eq 4) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-09 09:47:22"

foreach ($hl in $dede.Hyperlinks) {
    $r = $hl.Range
    if ($r.Row -eq 3 -and $r.Column -eq 4) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
